$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H3 becomes a new "Team meeting" note (replacing the old "Worked on week 2
# tutorial " text that used to live there) and picks up the shaded
# "Team meeting" style already used by G3/D3 (xf index 10: shaded fill +
# same border/wrap as the plain text cells).
$ws.Range("H3").Value = "Team meeting continuted. Guzzi has problems with his enum."
$ws.Range("G3").Copy()
$ws.Range("H3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# The "Worked on week 2 tutorial " text that used to sit in H3 is now
# duplicated into I3 and J3 (previously empty), keeping their original
# (unshaded) style.
$ws.Range("I3").Value = "Worked on week 2 tutorial "
$ws.Range("J3").Value = "Worked on week 2 tutorial "

# Update the visible selection/scroll position to match the saved view.
$ws.Range("J3").Select()
